# Apply cryptos.xlsx data refresh (prices / 1h volume %) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.910.44"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.888.14"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7743"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.89"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3103"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.66"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07171"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08595"
$ws.Range("E11").Value = "  +5.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7646"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "1.892.44"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.364"
$ws.Range("E14").Value = "  -2.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.80"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.154"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "29.930.67"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.77"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.17"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007820"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "2.207.65"
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9979"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.956"
$ws.Range("E23").Value = "  -2.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1641"
$ws.Range("E25").Value = "  +4.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.358"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.40"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.037"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.439"
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.533"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.108"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7468"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.696"
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.781"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4466"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "1.107.97"
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.083"
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "73.07"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.66"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.871"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.607"
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.116.91"
$ws.Range("E50").Value = "  +2.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.994"
$ws.Range("E51").Value = "  -1.00%  "
